$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F2 value
$ws.Range("F2").Value = 0.8874061718098415

# Add new row 3
$ws.Range("A3").Value = "gated_cnn"
$ws.Range("B3").Value = 0.001
$ws.Range("C3").Value = 128
$ws.Range("D3").Value = 64
$ws.Range("E3").Value = "max"
$ws.Range("F3").Value = 0.8807339449541285

# Add new row 4
$ws.Range("A4").Value = "gated_cnn"
$ws.Range("B4").Value = 0.001
$ws.Range("C4").Value = 128
$ws.Range("D4").Value = 256
$ws.Range("E4").Value = "avg"
$ws.Range("F4").Value = 0.8840700583819849
